$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3560
$ws.Range("I62").Value = 5166.6665
$ws.Range("J62").Value = 1150
$ws.Range("K62").Value = 5166.6665
$ws.Range("L62").Value = 1150
$ws.Range("M62").Value = -4542.6665
$ws.Range("N62").Value = -2398
$ws.Range("H65").Value = 3560
$ws.Range("I65").Value = 5166.6665
$ws.Range("J65").Value = 1150
$ws.Range("K65").Value = 25833.3325
$ws.Range("L65").Value = 5750
$ws.Range("M65").Value = -22713.3325
$ws.Range("N65").Value = -11990
$ws.Range("H138").Value = 2038.2245
$ws.Range("I138").Value = 1334.9259
$ws.Range("J138").Value = 2901.3635
$ws.Range("K138").Value = 4004.7777
$ws.Range("L138").Value = 8704.0905
$ws.Range("M138").Value = 1135.2223
$ws.Range("N138").Value = -18984.0905
$ws.Range("H139").Value = 95739
$ws.Range("J139").Value = 95739
$ws.Range("L139").Value = 95739
$ws.Range("N139").Value = -106019
$ws.Range("H140").Value = 55000
$ws.Range("J140").Value = 55000
$ws.Range("L140").Value = 55000
$ws.Range("N140").Value = -65360

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H70").Value = 48000
$ws.Range("J70").Value = 48000
$ws.Range("L70").Value = 48000
$ws.Range("N70").Value = -48540
$ws.Range("H73").Value = 48000
$ws.Range("J73").Value = 48000
$ws.Range("L73").Value = 48000
$ws.Range("N73").Value = -49872
$ws.Range("H122").Value = 9107.714
$ws.Range("I122").Value = 11028.182
$ws.Range("J122").Value = 2066
$ws.Range("K122").Value = 33084.546
$ws.Range("L122").Value = 6198
$ws.Range("M122").Value = -30634.546
$ws.Range("N122").Value = -11098

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1852.4166
$ws.Range("I94").Value = 1623.8
$ws.Range("J94").Value = 2015.7142
$ws.Range("K94").Value = 1623.8
$ws.Range("L94").Value = 2015.7142
$ws.Range("M94").Value = -1172.8
$ws.Range("N94").Value = -2917.7142

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1966.4138
$ws.Range("I58").Value = 807
$ws.Range("K58").Value = 807
$ws.Range("M58").Value = -604
$ws.Range("H132").Value = 22734430
$ws.Range("I132").Value = 35722692
$ws.Range("J132").Value = 4971
$ws.Range("K132").Value = 107168076
$ws.Range("L132").Value = 14913
$ws.Range("M132").Value = -107165546
$ws.Range("N132").Value = -19973
$ws.Range("H136").Value = 1966.4138
$ws.Range("I136").Value = 807
$ws.Range("K136").Value = 2421
$ws.Range("M136").Value = 129
$ws.Range("H140").Value = 32247.111
$ws.Range("J140").Value = 32247.111
$ws.Range("L140").Value = 32247.111
$ws.Range("N140").Value = -42607.111

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 65.56
$ws.Range("I2").Value = 113.545456
$ws.Range("J2").Value = 27.857143
$ws.Range("K2").Value = 681.272736
$ws.Range("L2").Value = 167.142858
$ws.Range("M2").Value = -568.272736
$ws.Range("N2").Value = -393.142858
$ws.Range("H38").Value = 236.66667
$ws.Range("I38").Value = 550
$ws.Range("J38").Value = 80
$ws.Range("K38").Value = 1650
$ws.Range("L38").Value = 240
$ws.Range("M38").Value = -1303
$ws.Range("N38").Value = -934
$ws.Range("H80").Value = 2090.5454
$ws.Range("J80").Value = 2375
$ws.Range("L80").Value = 7125
$ws.Range("N80").Value = -8997
$ws.Range("H83").Value = 2090.5454
$ws.Range("J83").Value = 2375
$ws.Range("L83").Value = 21375
$ws.Range("N83").Value = -30735

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 15000
$ws.Range("J52").Value = 15000
$ws.Range("L52").Value = 15000
$ws.Range("N52").Value = -15518
$ws.Range("H122").Value = 3032058.2
$ws.Range("I122").Value = 3510367
$ws.Range("J122").Value = 2769.3333
$ws.Range("K122").Value = 10531101
$ws.Range("L122").Value = 8307.999899999999
$ws.Range("M122").Value = -10528651
$ws.Range("N122").Value = -13207.9999
$ws.Range("H124").Value = 60780
$ws.Range("J124").Value = 60780
$ws.Range("L124").Value = 60780
$ws.Range("N124").Value = -70600
$ws.Range("H130").Value = 39950
$ws.Range("J130").Value = 39950
$ws.Range("L130").Value = 39950
$ws.Range("N130").Value = -49990
$ws.Range("H132").Value = 4771.6665
$ws.Range("I132").Value = 3432.75
$ws.Range("K132").Value = 10298.25
$ws.Range("M132").Value = -7768.25
$ws.Range("H138").Value = 58732.832
$ws.Range("J138").Value = 58732.832
$ws.Range("L138").Value = 58732.832
$ws.Range("N138").Value = -69012.83199999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 17247842
$ws.Range("I136").Value = 33336134
$ws.Range("J136").Value = 10386.786
$ws.Range("K136").Value = 100008402
$ws.Range("L136").Value = 31160.358
$ws.Range("M136").Value = -100005852
$ws.Range("N136").Value = -36260.358
$ws.Range("H139").Value = 46469.4
$ws.Range("J139").Value = 47116
$ws.Range("L139").Value = 47116
$ws.Range("N139").Value = -57396

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("H62").Value = 9720.888999999999
$ws.Range("I62").Value = 4511.5
$ws.Range("J62").Value = 13888.4
$ws.Range("K62").Value = 4511.5
$ws.Range("L62").Value = 13888.4
$ws.Range("M62").Value = -3887.5
$ws.Range("N62").Value = -15136.4
$ws.Range("H65").Value = 9720.888999999999
$ws.Range("I65").Value = 4511.5
$ws.Range("J65").Value = 13888.4
$ws.Range("K65").Value = 22557.5
$ws.Range("L65").Value = 69442
$ws.Range("M65").Value = -19437.5
$ws.Range("N65").Value = -75682
$ws.Range("H122").Value = 2330.7917
$ws.Range("I122").Value = 2309.75
$ws.Range("J122").Value = 2372.875
$ws.Range("K122").Value = 6929.25
$ws.Range("L122").Value = 7118.625
$ws.Range("M122").Value = -4479.25
$ws.Range("N122").Value = -12018.625
$ws.Range("H126").Value = 4002.6843
$ws.Range("I126").Value = 1459.3334
$ws.Range("J126").Value = 8362.714
$ws.Range("K126").Value = 4378.0002
$ws.Range("L126").Value = 25088.142
$ws.Range("M126").Value = -1908.0002
$ws.Range("N126").Value = -30028.142
$ws.Range("N57").ClearContents()
